$d = $word.ActiveDocument

$find = "a felhasználó által megadott adatokkal."
$replace = "a felhasználó által megadott adatokkal. Minden adat egy szóból állhat, egyik adat sem tárol el space karaktert."

$d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $replace, 2)
